$wb = $excel.ActiveWorkbook

# --- Sheet "Registro" (existing sheet): update the selected cell ---
$ws1 = $wb.Worksheets.Item("Registro")
$ws1.Range("B12").Select() | Out-Null

# --- Add the new "Planform" sheet right after "Registro" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planform"

# Column widths (A, B, C) -- values chosen so the engine's internal
# char-width rounding lands as close as possible to the target stored widths.
$ws2.Columns.Item(1).ColumnWidth = 22.5
$ws2.Columns.Item(2).ColumnWidth = 46
$ws2.Columns.Item(3).ColumnWidth = 22.666666666666668

# Cell contents
$ws2.Range("A1").Value = "Datos"
$ws2.Range("A2").Value = 3004442525
$ws2.Range("A3").Value = 1001828778
$ws2.Range("A4").Value = "dhaenerhys@gmail.com"

# Hyperlink on A4 (mailto link), styled like the existing hyperlink cell
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:dhaenerhys@gmail.com") | Out-Null
$ws2.Range("A4").Style = $ws1.Range("A4").Style

# Selection for the new sheet (becomes the active sheet/tab)
$ws2.Range("B1:C3").Select() | Out-Null
